$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 8.676671000000001
$ws.Range("N2").Value = 26.030013
$ws.Range("O2").Value = 0.1325240072999665
$ws.Range("P2").Value = 0.1325240072999665
$ws.Range("Q2").Value = 0.6183632043806667
$ws.Range("R2").Value = 5.565268839425999
$ws.Range("S2").Value = 0.1325240072999665
$ws.Range("T2").Value = 0.1325240072999665

# Row 3
$ws.Range("M3").Value = 37.74750533333334
$ws.Range("O3").Value = 0.5765403197090441
$ws.Range("P3").Value = 0.576540319709044
$ws.Range("Q3").Value = 2.690164045092445
$ws.Range("S3").Value = 0.5765403197090441
$ws.Range("T3").Value = 0.576540319709044

# Row 4
$ws.Range("M4").Value = 19.04827033333333
$ws.Range("N4").Value = 57.144811
$ws.Range("O4").Value = 0.2909356729909895
$ws.Range("P4").Value = 0.2909356729909895
$ws.Range("Q4").Value = 1.357519431269111
$ws.Range("R4").Value = 12.217674881422
$ws.Range("S4").Value = 0.2909356729909895
$ws.Range("T4").Value = 0.2909356729909895
